$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 262.95206511980251
$ws.Range("C2").Value = 289.18338212380695
$ws.Range("D2").Value = 259.9637642171167
$ws.Range("E2").Value = 294.60394532936505

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 259.5601453432763
$ws.Range("C3").Value = 293.65608461176885
$ws.Range("D3").Value = 259.22057767643787
$ws.Range("E3").Value = 300.88329427136955

# Update selection to B1:E3
$ws.Range("B1:E3").Select()
